$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert two new columns (D:E) for the new quarter-end periods (2018-12-31, 2018-09-30),
#    shifting all existing quarterly columns two positions to the right (old D -> F, ... old K -> M).
$ws.Columns("D:E").Insert()

# 2) Copy number formatting (date / thousands) from the (now-shifted) first data column F
#    into the two freshly inserted blank columns, per contiguous statement block.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)

$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)

$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 3) Write every updated cell value: the two new quarter columns (D,E) plus the handful of
#    restated figures in the pre-existing (now shifted) columns for this refresh.
$updates = @(
  @{r=7; c=4; v=43465},
  @{r=7; c=5; v=43373},
  @{r=8; c=4; v=1360600},
  @{r=8; c=5; v=1373300},
  @{r=9; c=4; v=1266800},
  @{r=9; c=5; v=1280000},
  @{r=10; c=4; v=93800},
  @{r=10; c=5; v=93300},
  @{r=12; c=4; v="NA"},
  @{r=12; c=5; v="NA"},
  @{r=13; c=4; v=0},
  @{r=13; c=5; v=0},
  @{r=14; c=4; v=307700},
  @{r=14; c=5; v=400},
  @{r=15; c=4; v=0},
  @{r=15; c=5; v=0},
  @{r=17; c=4; v=1654300},
  @{r=17; c=5; v=1363700},
  @{r=18; c=4; v=-293700},
  @{r=18; c=5; v=9600},
  @{r=20; c=4; v=600},
  @{r=20; c=5; v=600},
  @{r=21; c=4; v=-268600},
  @{r=21; c=5; v=34600},
  @{r=22; c=4; v=10700},
  @{r=22; c=5; v=10200},
  @{r=23; c=4; v=-303800},
  @{r=23; c=5; v=0},
  @{r=24; c=4; v=-4700},
  @{r=24; c=5; v=-100},
  @{r=25; c=4; v=0},
  @{r=25; c=5; v=0},
  @{r=26; c=4; v=-299100},
  @{r=26; c=5; v=200},
  @{r=27; c=4; v=-299100},
  @{r=27; c=5; v=200},
  @{r=28; c=4; v=0},
  @{r=28; c=5; v=0},
  @{r=29; c=4; v=1100},
  @{r=29; c=5; v="NA"},
  @{r=30; c=4; v=0},
  @{r=30; c=5; v=0},
  @{r=31; c=4; v=0},
  @{r=31; c=5; v=0},
  @{r=32; c=4; v=-600},
  @{r=32; c=5; v=-600},
  @{r=33; c=4; v=-298000},
  @{r=33; c=5; v=200},
  @{r=34; c=4; v=0},
  @{r=34; c=5; v=0},
  @{r=35; c=4; v=-298000},
  @{r=35; c=5; v=200},
  @{r=38; c=4; v=43465},
  @{r=38; c=5; v=43373},
  @{r=41; c=4; v=9500},
  @{r=41; c=5; v=8200},
  @{r=42; c=4; v=0},
  @{r=42; c=5; v=0},
  @{r=43; c=4; v=326600},
  @{r=43; c=5; v=358200},
  @{r=44; c=4; v=210600},
  @{r=44; c=5; v=169900},
  @{r=45; c=4; v=9600},
  @{r=45; c=5; v=13500},
  @{r=46; c=4; v=556300},
  @{r=46; c=5; v=549700},
  @{r=47; c=4; v=0},
  @{r=47; c=5; v=0},
  @{r=48; c=4; v=34500},
  @{r=48; c=5; v=40900},
  @{r=49; c=4; v=880900},
  @{r=49; c=5; v=1205100},
  @{r=50; c=4; v=0},
  @{r=50; c=5; v=0},
  @{r=51; c=4; v=0},
  @{r=51; c=5; v=0},
  @{r=52; c=4; v=4700},
  @{r=52; c=5; v=4900},
  @{r=53; c=4; v=0},
  @{r=53; c=5; v=0},
  @{r=54; c=4; v=1476400},
  @{r=54; c=5; v=1800700},
  @{r=57; c=4; v=308100},
  @{r=57; c=5; v=315300},
  @{r=58; c=4; v=187800},
  @{r=58; c=5; v=189800},
  @{r=59; c=4; v=62700},
  @{r=59; c=5; v=73500},
  @{r=60; c=4; v=558600},
  @{r=60; c=5; v=578500},
  @{r=61; c=4; v=438400},
  @{r=61; c=5; v=440600},
  @{r=62; c=4; v=14300},
  @{r=62; c=5; v=16800},
  @{r=63; c=4; v=0},
  @{r=63; c=5; v=0},
  @{r=64; c=4; v=0},
  @{r=64; c=5; v=0},
  @{r=65; c=4; v=0},
  @{r=65; c=5; v=0},
  @{r=66; c=4; v=1011300},
  @{r=66; c=5; v=1035800},
  @{r=68; c=4; v=0},
  @{r=68; c=5; v=0},
  @{r=69; c=4; v=0},
  @{r=69; c=5; v=0},
  @{r=70; c=4; v=0},
  @{r=70; c=5; v=0},
  @{r=71; c=4; v=0},
  @{r=71; c=5; v=0},
  @{r=72; c=4; v=-210600},
  @{r=72; c=5; v=87400},
  @{r=73; c=4; v=0},
  @{r=73; c=5; v=0},
  @{r=74; c=4; v=0},
  @{r=74; c=5; v=0},
  @{r=75; c=4; v=0},
  @{r=75; c=5; v=0},
  @{r=76; c=4; v=465100},
  @{r=76; c=5; v=764900},
  @{r=77; c=4; v=0},
  @{r=77; c=5; v=0},
  @{r=80; c=4; v=43465},
  @{r=80; c=5; v=43373},
  @{r=81; c=4; v=-298000},
  @{r=81; c=5; v=200},
  @{r=83; c=4; v=24600},
  @{r=83; c=5; v=24400},
  @{r=84; c=4; v=0},
  @{r=84; c=5; v=0},
  @{r=85; c=4; v=0},
  @{r=85; c=5; v=0},
  @{r=86; c=4; v=0},
  @{r=86; c=5; v=0},
  @{r=87; c=4; v=0},
  @{r=87; c=5; v=0},
  @{r=88; c=4; v=0},
  @{r=88; c=5; v=0},
  @{r=89; c=4; v=1800},
  @{r=89; c=5; v=-33400},
  @{r=91; c=4; v=6700},
  @{r=91; c=5; v=-2400},
  @{r=92; c=4; v=0},
  @{r=92; c=5; v=0},
  @{r=93; c=4; v=0},
  @{r=93; c=5; v=0},
  @{r=94; c=4; v=6500},
  @{r=94; c=5; v=-5100},
  @{r=96; c=4; v=0},
  @{r=96; c=5; v=0},
  @{r=97; c=4; v=0},
  @{r=97; c=5; v=0},
  @{r=98; c=4; v=0},
  @{r=98; c=5; v=0},
  @{r=99; c=4; v=0},
  @{r=99; c=5; v=0},
  @{r=100; c=4; v=-7000},
  @{r=100; c=5; v=39400},
  @{r=101; c=4; v=0},
  @{r=101; c=5; v=0},
  @{r=102; c=4; v=1300},
  @{r=102; c=5; v=900},
  @{r=9; c=8; v=1082500},
  @{r=9; c=9; v=1059900},
  @{r=10; c=8; v=72600},
  @{r=10; c=9; v=65100},
  @{r=14; c=6; v=100},
  @{r=14; c=7; v=900},
  @{r=14; c=8; v=1800},
  @{r=14; c=9; v=400},
  @{r=14; c=10; v=300},
  @{r=47; c=6; v=0},
  @{r=47; c=7; v=0},
  @{r=47; c=8; v=0},
  @{r=47; c=9; v=0},
  @{r=47; c=10; v=0},
  @{r=48; c=8; v=78000},
  @{r=57; c=8; v=384700},
  @{r=59; c=8; v=67000},
  @{r=91; c=6; v=-3200},
  @{r=91; c=7; v=-2300},
  @{r=91; c=8; v=-6600},
  @{r=91; c=9; v=0},
  @{r=91; c=10; v=0}
)

foreach ($u in $updates) {
  $ws.Cells.Item($u.r, $u.c).Value = $u.v
}

